$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column contains text values that look numeric (e.g. "1.003"),
# so force text format before assigning, then restore default style
# so the cell style index matches the original (unstyled) cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.013.19"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "1.890.19"
$ws.Range("E3").Value = "  -3.85%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.83%  "
$ws.Range("D5").Value = "325.87"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").Value = "0.4591"
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "51.45"
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").Value = "0.08229"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("D11").Value = "1.036"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "21.57"
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("D13").Value = "1.937.62"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "7.311"
$ws.Range("E14").Value = "  -4.68%  "
$ws.Range("D15").Value = "5.974"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "89.17"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "0.00001058"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "0.06579"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "17.57"
$ws.Range("E20").Value = "  -6.12%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "5.642"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").Value = "27.979.52"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "11.09"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "2.100.16"
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("D27").Value = "154.08"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "2.101"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "5.659"
$ws.Range("E30").Value = "  -5.25%  "
$ws.Range("D31").Value = "124.08"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "0.09540"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "0.9581"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("D34").Value = "1.461"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "3.630"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").Value = "5.464"
$ws.Range("E36").Value = "  -3.86%  "
$ws.Range("D37").Value = "0.02282"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "1.252"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "8.636"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "0.06104"
$ws.Range("D41").Value = "0.6089"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "1.307"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").Value = "0.5812"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "12.66"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "1.989"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "0.06882"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "110.19"
$ws.Range("E51").Value = "  -0.76%  "

# Restore default (unstyled) cell style for the D column so the
# saved style index matches the original workbook.
$ws.Range("D2:D51").Style = "Normal"
